$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(467).Insert()

$ws.Range("A467").Value = 8
$ws.Range("B467").Value = "Terminal La Palmera de La Serena"
$ws.Range("C467").Value = "Coquimbo"
$ws.Range("D467").Value = 45211
$ws.Range("E467").Value = 4
$ws.Range("F467").Value = 100112021
$ws.Range("G467").Value = "Ají"
$ws.Range("H467").Value = "Inferno"
$ws.Range("I467").Value = "Primera"
$ws.Range("J467").Value = 460
$ws.Range("K467").Value = 29000
$ws.Range("L467").Value = 30000
$ws.Range("M467").Value = 29500
$ws.Range("N467").Value = "$/caja 10 kilos"
$ws.Range("O467").Value = "Región de Arica y Parinacota"
$ws.Range("P467").Value = 2950
$ws.Range("Q467").Value = 10
$ws.Range("R467").Value = "Hortaliza"
